$wb = $excel.ActiveWorkbook

# --- T3 sheet: remove a few "o" (possible substrate inhibition) markers ---
$ws3 = $wb.Worksheets.Item("T3")
$ws3.Range("G2").ClearContents()
$ws3.Range("H6").ClearContents()
$ws3.Range("D8").ClearContents()
$ws3.Range("H10").ClearContents()
$ws3.Range("D10").Select()

# --- T5 sheet: adjust error markers now that PPO (oxidase) parameters were obtained ---
$ws5 = $wb.Worksheets.Item("T5")
$ws5.Range("H6").ClearContents()
$ws5.Range("B7").ClearContents()
$ws5.Range("G7").Value = "o"
$ws5.Range("G10").Value = "o"
$ws5.Range("G12").Value = "o"
$ws5.Range("H13").Value = "o"
$ws5.Range("E15").Value = "o"
$ws5.Range("F15").Select()

# --- T6 sheet: becomes the active/selected sheet ---
$ws6 = $wb.Worksheets.Item("T6")
$ws6.Activate()
$ws6.Range("F19").Select()
